# The workbook originally has an unused/blank column E (between the
# "notes " header in D and the SiO2 data starting in F) and an unused/blank
# column Q (between the K2O data in P and the "tot" column in R).
# Both were removed from the sheet, shifting all the analytical data
# columns one (and the final "tot" column two) positions to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank spacer columns (E first, then what was originally Q,
# which is now column P after the first deletion shifts everything left).
$ws.Columns("E").Delete()
$ws.Columns("P").Delete()

# Reflect the final selection left behind on the worksheet: the whole of
# the new last column (P), which used to be column R/"tot".
$ws.Range("P1:P1048576").Select()
